# weekly update, 2021-11-17, SARS-CoV-2
$wb = $excel.ActiveWorkbook

# 1) Rename sheet 5
$sheetLinker = $wb.Worksheets.Item(5)
$sheetLinker.Name = "nsp3_pred_linker-Ecto-TM2"

# 2) Insert new PDB entry row into nsp3_PL2pro (sheet 2), after "7ofs" row (row 20), before "7sgu" row (row 21)
$wsPL2pro = $wb.Worksheets.Item(2)
$wsPL2pro.Rows.Item(21).Insert()
$wsPL2pro.Cells.Item(21,1).Value = "7sqe"
$wsPL2pro.Cells.Item(21,2).Value = 2
$wsPL2pro.Cells.Item(21,3).Value = "PAPAIN-LIKE PROTEASE OF SARS COV-2, C111S MUTANT, IN COMPLEX WITH JUN9-84-3 INHIBITOR"
$wsPL2pro.Cells.Item(21,4).Value = "X-RAY DIFFRACTION"
$wsPL2pro.Cells.Item(21,5).Value = "2021-11-05"

# 3) Insert new PDB id into Organisms sheet (sheet 8), in the
#    "severe acute respiratory syndrome coronavirus2" block, right after the
#    row holding "7ofs" (row 146), before the row holding "5rua" (row 150)
$wsOrganisms = $wb.Worksheets.Item(8)
$wsOrganisms.Rows.Item(150).Insert()
$wsOrganisms.Cells.Item(150,2).Value = "7sqe"
